# "revamp macrobuilder use of template, get dossier class working"
#
# The three per-sample worksheets ("set1"/"set2"/"set3") are renamed to
# the new "grid" naming used by the revamped template ("grid1"/"grid2"/
# "grid3"), and the workbook's active/selected tab moves from the first
# sheet to the third ("grid3") -- i.e. activeTab goes from 0 to 2 and the
# tabSelected flag moves off grid1 and onto grid3.

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("set1").Name = "grid1"
$wb.Worksheets.Item("set2").Name = "grid2"
$wb.Worksheets.Item("set3").Name = "grid3"

$wb.Worksheets.Item("grid3").Activate()
